$d = $word.ActiveDocument

# Remove all existing content (leaves a single empty paragraph behind)
$d.Content.Delete()

$para1 = '    To efficiently meet my goal of visiting 128 cities across the US (and adjacent areas of Canada), I have divided those cities into travel routes.  This visualization scores the weather along each route for each possible two-week travel window.  For this project, pleasant weather means temperatures between 55 and 75 degrees and less than 0.2 inches of precipitation per hour.  I only consider weather between 7 AM and 10 PM.'

$para2 = '    Each chart scores the weather for each two-week travel window (x-axis) along each route (lines).  The x-axis indicates the first day of the two-week window.  The score (y-axes) measures the probability that the weather will be pleasant in each city on the route at each hour of the day.  The maximum score is 1.0, and the minimum is 0. See github.com/sjoshuam/us_travels for more on each route.'

# Fill the first (remaining) paragraph with the new explanation text
$r1 = $d.Paragraphs(1).Range
$r1.InsertAfter($para1)

# Create a second paragraph and fill it with the supporting text
$rEnd = $d.Content
$rEnd.Collapse(0)
$rEnd.InsertParagraphAfter()

$r2 = $d.Paragraphs(2).Range
$r2.InsertAfter($para2)
